# Round 5 of FDA review.
# The "submission_date" field is renamed to "date_submitted" and its
# description gains a trailing note about the date format. This also
# shifts a couple of shared-string-driven rebuilds (row 2 picks up an
# explicit row height, and the active selection moves to D3 where the
# edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: field name + description for the "date received by CTP" field.
$ws.Range("B3").Value = "date_submitted"
$ws.Range("D3").Value = "Date report was received by CTP; this is the earliest date of report receipt, either to Safety Reporting Portal (SRP) or by other means. Date follows format: ``YYYYmmdd``."

# Row 2 ends up with an explicit row height of 17 after the edit.
$ws.Rows.Item(2).RowHeight = 17

# Leave the selection on the cell that was actually edited.
$ws.Range("D3").Select() | Out-Null
